# Insert a new weekly price record at row 188 (2023-01-20 / Albahaca),
# pushing the existing rows 188-218 down to 189-219.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(188).Insert()

$ws.Range("A188").Value = 3
$ws.Range("B188").Value = "Femacal de La Calera"
$ws.Range("C188").Value = "Coquimbo"
$ws.Range("D188").Value = 44946
$ws.Range("E188").Value = 5
$ws.Range("F188").Value = 100112052
$ws.Range("G188").Value = "Albahaca"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 115
$ws.Range("K188").Value = 4000
$ws.Range("L188").Value = 4500
$ws.Range("M188").Value = 4239
$ws.Range("N188").Value = "$/docena de matas"
$ws.Range("O188").Value = "Provincia de Quillota"
$ws.Range("P188").Value = 706
$ws.Range("Q188").Value = 6
$ws.Range("R188").Value = "Hortaliza"
